$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-table statistics (ME, MAE, MSE, RMSE, SE, N) for rows 2-10
# (Q1..Q9), reflecting an additional preprocessed ifo GDP component
# observation folded into the ifoCAST matched-error calculations.

$ws.Cells.Item(2, 2).Value = -0.1679391225927079
$ws.Cells.Item(2, 3).Value = 0.6512668612981908
$ws.Cells.Item(2, 4).Value = 1.077333483753495
$ws.Cells.Item(2, 5).Value = 1.037946763448634
$ws.Cells.Item(2, 6).Value = 1.047290627843072
$ws.Cells.Item(2, 7).Value = 23

$ws.Cells.Item(3, 2).Value = 0.6108509255840642
$ws.Cells.Item(3, 3).Value = 0.8713162755896118
$ws.Cells.Item(3, 4).Value = 1.94510215081865
$ws.Cells.Item(3, 5).Value = 1.394669190460107
$ws.Cells.Item(3, 6).Value = 1.283284338193333
$ws.Cells.Item(3, 7).Value = 22

$ws.Cells.Item(4, 2).Value = 0.5135409635362268
$ws.Cells.Item(4, 3).Value = 1.149401313682296
$ws.Cells.Item(4, 4).Value = 3.430862853712828
$ws.Cells.Item(4, 5).Value = 1.852258851703192
$ws.Cells.Item(4, 6).Value = 1.8235941048126
$ws.Cells.Item(4, 7).Value = 21

$ws.Cells.Item(5, 2).Value = 0.6385223778103771
$ws.Cells.Item(5, 3).Value = 0.8102001327720327
$ws.Cells.Item(5, 4).Value = 1.083726471416818
$ws.Cells.Item(5, 5).Value = 1.041021840028737
$ws.Cells.Item(5, 6).Value = 0.8435611508437559
$ws.Cells.Item(5, 7).Value = 20

$ws.Cells.Item(6, 2).Value = 0.5184854322233537
$ws.Cells.Item(6, 3).Value = 0.6702212091134594
$ws.Cells.Item(6, 4).Value = 0.6945837070328604
$ws.Cells.Item(6, 5).Value = 0.8334168866976841
$ws.Cells.Item(6, 6).Value = 0.6703802697182577
$ws.Cells.Item(6, 7).Value = 19

$ws.Cells.Item(7, 2).Value = 0.3259858714718396
$ws.Cells.Item(7, 3).Value = 0.5159630503235952
$ws.Cells.Item(7, 4).Value = 0.3780268743366828
$ws.Cells.Item(7, 5).Value = 0.6148389011250693
$ws.Cells.Item(7, 6).Value = 0.5364195870263428
$ws.Cells.Item(7, 7).Value = 18

$ws.Cells.Item(8, 2).Value = 0.345936562473689
$ws.Cells.Item(8, 3).Value = 0.3960995140440492
$ws.Cells.Item(8, 4).Value = 0.2072861800771614
$ws.Cells.Item(8, 5).Value = 0.45528692060849
$ws.Cells.Item(8, 6).Value = 0.3051064642012103
$ws.Cells.Item(8, 7).Value = 17

$ws.Cells.Item(9, 2).Value = 0.3283187899062386
$ws.Cells.Item(9, 3).Value = 0.3609349132557007
$ws.Cells.Item(9, 4).Value = 0.1816422354358933
$ws.Cells.Item(9, 5).Value = 0.4261950673528418
$ws.Cells.Item(9, 6).Value = 0.2838354343252694
$ws.Cells.Item(9, 7).Value = 12

$ws.Cells.Item(10, 2).Value = 0.3306938847573825
$ws.Cells.Item(10, 3).Value = 0.3306938847573825
$ws.Cells.Item(10, 4).Value = 0.1666290365104273
$ws.Cells.Item(10, 5).Value = 0.4082022005212947
$ws.Cells.Item(10, 6).Value = 0.2584873103466553
$ws.Cells.Item(10, 7).Value = 7
